$wb = $excel.ActiveWorkbook

# --- Current layout before edit ---
#   Sheet 1: "总计"      (summary)
#   Sheet 2: "2022-Q3"   (fund-detail data, currently for Q3)
#
# --- Target layout after edit ---
#   Sheet 1: "总计"      (summary, now with an extra Q4 row)
#   Sheet 2: "2022-Q4"   (fund-detail data, new quarter)
#   Sheet 3: "2022-Q3"   (fund-detail data, same content as before, just relocated)

$wsTotal = $wb.Worksheets.Item(1)
$wsQ3Old = $wb.Worksheets.Item(2)

# 1) Rename the existing fund-detail sheet to "2022-Q4" - it keeps its sheetId,
#    we will overwrite its contents with the new quarter's data below.
$wsQ3Old.Name = "2022-Q4"
$wsQ4 = $wsQ3Old

# 2) Insert a brand-new sheet right after it and name it "2022-Q3"; this new
#    sheet will receive the original Q3 fund-detail content (unchanged).
$wsQ3New = $wb.Worksheets.Add($null, $wsQ4)
$wsQ3New.Name = "2022-Q3"

# 3) Relocate the original Q3 fund-detail content (still sitting in $wsQ4 under
#    its old name) onto the new "2022-Q3" sheet, preserving values + styles.
#    (Column A is copied separately, starting at row 2, to avoid materialising
#    a spurious blank A1 cell that a plain A1:H3 copy would introduce.)
$wsQ4.Range("B1:H3").Copy($wsQ3New.Range("B1:H3"))
$wsQ4.Range("A2:A3").Copy($wsQ3New.Range("A2:A3"))
$wsQ3New.Activate()

# 4) Clear out the old content on $wsQ4 and write the new 2022-Q4 fund data.
$wsQ4.Range("A1:H3").ClearContents()

# Columns that must stay TEXT (not get auto-coerced to numbers), matching the
# original data's inline-string typing (fund codes, percentages, scaled amounts).
$wsQ4.Range("B1:G4").NumberFormat = "@"

$wsQ4.Cells.Item(1,2).Value = "基金代码"
$wsQ4.Cells.Item(1,3).Value = "基金名称"
$wsQ4.Cells.Item(1,4).Value = "基金规模"
$wsQ4.Cells.Item(1,5).Value = "股票总仓位"
$wsQ4.Cells.Item(1,6).Value = "仓位占比"
$wsQ4.Cells.Item(1,7).Value = "持有市值(亿元)"
$wsQ4.Cells.Item(1,8).Value = "仓位排名"

$wsQ4.Cells.Item(2,1).Value = 0
$wsQ4.Cells.Item(2,2).Value = "003567"
$wsQ4.Cells.Item(2,3).Value = "华夏行业景气混合"
$wsQ4.Cells.Item(2,4).Value = "109.60"
$wsQ4.Cells.Item(2,5).Value = "93.65"
$wsQ4.Cells.Item(2,6).Value = "2.01"
$wsQ4.Cells.Item(2,7).Value = "2.2030"
$wsQ4.Cells.Item(2,8).Value = 8

$wsQ4.Cells.Item(3,1).Value = 1
$wsQ4.Cells.Item(3,2).Value = "016250"
$wsQ4.Cells.Item(3,3).Value = "华夏远见成长一年持有混合A"
$wsQ4.Cells.Item(3,4).Value = "9.60"
$wsQ4.Cells.Item(3,5).Value = "88.62"
$wsQ4.Cells.Item(3,6).Value = "5.45"
$wsQ4.Cells.Item(3,7).Value = "0.5232"
$wsQ4.Cells.Item(3,8).Value = 3

$wsQ4.Cells.Item(4,1).Value = 2
$wsQ4.Cells.Item(4,2).Value = "016251"
$wsQ4.Cells.Item(4,3).Value = "华夏远见成长一年持有混合C"
$wsQ4.Cells.Item(4,4).Value = "2.97"
$wsQ4.Cells.Item(4,5).Value = "88.62"
$wsQ4.Cells.Item(4,6).Value = "5.45"
$wsQ4.Cells.Item(4,7).Value = "0.1619"
$wsQ4.Cells.Item(4,8).Value = 3

# Re-apply the proper header/index-column formatting (bold, centered, bordered
# "s=2" style used by the 总计 sheet's header row) without touching the text
# values just written - PasteSpecial(xlPasteFormats) copies format only.
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A4").PasteSpecial(-4122)

# The data rows (2-4) should carry the plain default style (as in the source
# data) rather than the leftover border/bold/text-numfmt formatting from the
# sheet's previous life as the Q3 data sheet plus our temporary "@" format.
$wsTotal.Range("C2").Copy()
$wsQ4.Range("B2:G4").PasteSpecial(-4122)

# 5) Update the "总计" summary sheet: add a new Q4 row above the existing Q3
#    row, shifting Q3's values down and renumbering its index.
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Cells.Item(3,1).Value = 1
$wsTotal.Cells.Item(3,2).Value = "2022-Q3"
$wsTotal.Cells.Item(3,3).Value = 2
$wsTotal.Cells.Item(3,4).Value = 0.02

$wsTotal.Cells.Item(2,2).Value = "2022-Q4"
$wsTotal.Cells.Item(2,3).Value = 3
$wsTotal.Cells.Item(2,4).Value = 2.89
